$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is numeric-looking need an explicit Text format,
# otherwise Excel auto-converts the string into a number (losing the exact text).
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '68.781.65'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '2.511.39'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '592.18'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').Value = '174.53'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').Value = '2.510.25'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '0.150'
$ws.Range('E10').Value = '  +4.40%  '
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('D12').Value = '5.01'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').Value = '0.335'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').Value = '2.957.22'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').Value = '25.74'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '68.580.67'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').Value = '0.0000172'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '2.507.04'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '363.95'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('D20').Value = '7.55'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').Value = '10.88'
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('D22').Value = '4.01'
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '70.17'
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('D25').Value = '4.16'
$ws.Range('E25').Value = '  -3.48%  '
$ws.Range('D26').Value = '8.92'
$ws.Range('E26').Value = '  -3.65%  '
$ws.Range('E27').Value = '  -7.45%  '
$ws.Range('D28').Value = '2.635.63'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '510.58'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').Value = '0.0₃0878'
$ws.Range('E31').Value = '  -4.57%  '
$ws.Range('D32').Value = '7.72'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('E33').Value = '  -4.00%  '
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '161.66'
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').Value = '0.118'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('D38').Value = '18.51'
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  -0.33%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '1.31'
$ws.Range('E40').Value = '  -2.28%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '1.70'
$ws.Range('E42').Value = '  -3.12%  '
$ws.Range('D43').Value = '4.75'
$ws.Range('E43').Value = '  -3.28%  '
$ws.Range('D44').Value = '0.317'
$ws.Range('E44').Value = '  -4.37%  '
$ws.Range('D45').Value = '2.32'
$ws.Range('E45').Value = '  -5.14%  '
$ws.Range('D46').Value = '150.78'
$ws.Range('E46').Value = '  +3.38%  '
$ws.Range('D47').Value = '3.55'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').Value = '0.514'
$ws.Range('E48').Value = '  -0.95%  '
$ws.Range('D49').Value = '0.0737'
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('E50').Value = '  -5.47%  '
$ws.Range('D51').Value = '1.56'
$ws.Range('E51').Value = '  -2.25%  '
